$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = "rostralanteriorcingulate_1"
$ws.Range("E4").Value = -0.61
$ws.Range("F4").Value = -0.51

# Row 5
$ws.Range("B5").Value = "parahippocampal_1"
$ws.Range("E5").Value = -0.91
$ws.Range("F5").Value = ""

# Row 6
$ws.Range("B6").Value = "insula_1"
$ws.Range("E6").Value = -0.03
$ws.Range("F6").Value = ""

# Row 7
$ws.Range("B7").Value = "Right-Amygdala"
$ws.Range("E7").Value = -2.34
$ws.Range("F7").Value = -2.25

# Row 8
$ws.Range("B8").Value = "Right-Hippocampus"
$ws.Range("E8").Value = -0.93
$ws.Range("F8").Value = -0.66

# Row 9
$ws.Range("B9").Value = "lateralorbitofrontal_1"
$ws.Range("E9").Value = -0.64
$ws.Range("F9").Value = ""

# Row 10
$ws.Range("B10").Value = "parstriangularis_1"
$ws.Range("C10").Value = -2.12
$ws.Range("D10").Value = -2.36
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""

# Row 11
$ws.Range("B11").Value = "parsopercularis_1"
$ws.Range("C11").Value = -1.28
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""

# Row 12
$ws.Range("B12").Value = "parahippocampal_1"
$ws.Range("C12").Value = -1.59
$ws.Range("D12").Value = -1.91
$ws.Range("E12").Value = -0.88
$ws.Range("F12").Value = -0.91

# Row 13
$ws.Range("B13").Value = "Left-Amygdala"
$ws.Range("C13").Value = -2.11
$ws.Range("D13").Value = -2.15
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""

# Row 14
$ws.Range("B14").Value = "Left-Hippocampus"
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = -0.2
$ws.Range("F14").Value = -0.05

# Remove rows 15 and 16 (now stale) entirely
$ws.Range("A15:F16").ClearContents()
$ws.Rows("15:16").Delete()
